$d = $word.ActiveDocument

# Update the title date line.
$d.Content.Find.Execute("2023-12-17 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-12-18 Monday", 2)

# Update the multiplication table values. Cells are addressed positionally
# (table row/column) rather than via text Find/Replace because one new
# value ("98×68=6664") collides with an old value elsewhere in the table,
# which would make a blind global replace ambiguous/unsafe.
$tbl = $d.Tables.Item(1)

$newValues = @(
    @(1,  1, "57×40=2280"),
    @(1,  2, "76×98=7448"),
    @(1,  3, "25×98=2450"),
    @(1,  4, "26×95=2470"),
    @(1,  5, "82×55=4510"),

    @(5,  1, "40×14=560"),
    @(5,  2, "80×54=4320"),
    @(5,  3, "15×48=720"),
    @(5,  4, "72×79=5688"),
    @(5,  5, "78×88=6864"),

    @(10, 1, "64×84=5376"),
    @(10, 2, "80×46=3680"),
    @(10, 3, "75×55=4125"),
    @(10, 4, "78×59=4602"),
    @(10, 5, "76×43=3268"),

    @(15, 1, "95×32=3040"),
    @(15, 2, "11×63=693"),
    @(15, 3, "78×42=3276"),
    @(15, 4, "84×47=3948"),
    @(15, 5, "76×11=836"),

    @(20, 1, "92×44=4048"),
    @(20, 2, "13×51=663"),
    @(20, 3, "98×68=6664"),
    @(20, 4, "99×76=7524"),
    @(20, 5, "95×90=8550")
)

foreach ($entry in $newValues) {
    $row = $entry[0]
    $col = $entry[1]
    $value = $entry[2]
    $tbl.Cell($row, $col).Range.Text = $value
}

Write-Host "Applied date + table value updates"
